$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.46 = 30074.63 pesos`n✅ 30074.63 pesos = 7.43 = 954.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells N10/O10/N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 134
$ws2.Range("O10").Value = 4030
$ws2.Range("N12").Value = 4050
$ws2.Range("O12").Value = 128.52

Write-Host "edits applied"
